# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for
# rows 2-51. Values must remain stored as literal text (the source data
# uses locale-formatted numbers like "30.329.63" which are not valid
# Excel numerics, and percentage strings padded with spaces), so we force
# text entry without leaving a residual style/number-format on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.Style = "Normal"
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.329.63"
Set-TextValue $ws.Range("E2") "  -2.90%  "
Set-TextValue $ws.Range("D3") "1.936.88"
Set-TextValue $ws.Range("E3") "  -3.05%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "250.75"
Set-TextValue $ws.Range("E5") "  -2.34%  "
Set-TextValue $ws.Range("D6") "0.7245"
Set-TextValue $ws.Range("E6") "  -6.72%  "
Set-TextValue $ws.Range("E7") "  +0.02%  "
Set-TextValue $ws.Range("D8") "0.3317"
Set-TextValue $ws.Range("E8") "  -4.35%  "
Set-TextValue $ws.Range("D9") "27.79"
Set-TextValue $ws.Range("E9") "  -1.26%  "
Set-TextValue $ws.Range("D10") "0.07292"
Set-TextValue $ws.Range("E10") "  +0.49%  "
Set-TextValue $ws.Range("D11") "0.8090"
Set-TextValue $ws.Range("E11") "  -4.15%  "
Set-TextValue $ws.Range("D12") "0.08085"
Set-TextValue $ws.Range("E12") "  -1.35%  "
Set-TextValue $ws.Range("D13") "1.931.31"
Set-TextValue $ws.Range("E13") "  -3.46%  "
Set-TextValue $ws.Range("D14") "5.495"
Set-TextValue $ws.Range("E14") "  -2.63%  "
Set-TextValue $ws.Range("D15") "94.53"
Set-TextValue $ws.Range("E15") "  -6.31%  "
Set-TextValue $ws.Range("D16") "15.13"
Set-TextValue $ws.Range("E16") "  -3.31%  "
Set-TextValue $ws.Range("D17") "30.319.53"
Set-TextValue $ws.Range("E17") "  -2.92%  "
Set-TextValue $ws.Range("D18") "0.000008298"
Set-TextValue $ws.Range("E18") "  -0.50%  "
Set-TextValue $ws.Range("D19") "251.49"
Set-TextValue $ws.Range("E19") "  -7.75%  "
Set-TextValue $ws.Range("D20") "5.860"
Set-TextValue $ws.Range("E20") "  -2.31%  "
Set-TextValue $ws.Range("D21") "2.189.01"
Set-TextValue $ws.Range("E21") "  -2.90%  "
Set-TextValue $ws.Range("E22") "  +0.05%  "
Set-TextValue $ws.Range("D23") "1.001"
Set-TextValue $ws.Range("E23") "  +0.04%  "
Set-TextValue $ws.Range("D24") "6.982"
Set-TextValue $ws.Range("E24") "  -2.00%  "
Set-TextValue $ws.Range("D25") "9.767"
Set-TextValue $ws.Range("E25") "  -3.30%  "
Set-TextValue $ws.Range("D26") "164.12"
Set-TextValue $ws.Range("E26") "  -0.28%  "
Set-TextValue $ws.Range("D27") "2.377"
Set-TextValue $ws.Range("E27") "  -1.40%  "
Set-TextValue $ws.Range("D28") "19.29"
Set-TextValue $ws.Range("E28") "  -3.55%  "
Set-TextValue $ws.Range("D29") "0.1328"
Set-TextValue $ws.Range("E29") "  -6.04%  "
Set-TextValue $ws.Range("D30") "1.567"
Set-TextValue $ws.Range("E30") "  -2.64%  "
Set-TextValue $ws.Range("D31") "1.349"
Set-TextValue $ws.Range("E31") "  -1.42%  "
Set-TextValue $ws.Range("D32") "4.422"
Set-TextValue $ws.Range("E32") "  -5.10%  "
Set-TextValue $ws.Range("D33") "4.182"
Set-TextValue $ws.Range("E33") "  -6.62%  "
Set-TextValue $ws.Range("D34") "0.05205"
Set-TextValue $ws.Range("E34") "  -3.26%  "
Set-TextValue $ws.Range("D35") "1.281"
Set-TextValue $ws.Range("E35") "  +1.46%  "
Set-TextValue $ws.Range("D36") "0.7490"
Set-TextValue $ws.Range("E36") "  -4.99%  "
Set-TextValue $ws.Range("D37") "2.744"
Set-TextValue $ws.Range("E37") "  -1.30%  "
Set-TextValue $ws.Range("E38") "  -2.02%  "
Set-TextValue $ws.Range("D39") "2.822"
Set-TextValue $ws.Range("E39") "  -4.19%  "
Set-TextValue $ws.Range("D40") "78.86"
Set-TextValue $ws.Range("E40") "  -8.12%  "
Set-TextValue $ws.Range("D41") "6.374"
Set-TextValue $ws.Range("E41") "  -6.91%  "
Set-TextValue $ws.Range("D42") "0.4535"
Set-TextValue $ws.Range("E42") "  -3.44%  "
Set-TextValue $ws.Range("D43") "2.024"
Set-TextValue $ws.Range("E43") "  -5.56%  "
Set-TextValue $ws.Range("D44") "0.8484"
Set-TextValue $ws.Range("E44") "  -1.62%  "
Set-TextValue $ws.Range("D45") "1.000"
Set-TextValue $ws.Range("D46") "101.70"
Set-TextValue $ws.Range("E46") "  -3.29%  "
Set-TextValue $ws.Range("D47") "9.717"
Set-TextValue $ws.Range("E47") "  -5.18%  "
Set-TextValue $ws.Range("D48") "7.438"
Set-TextValue $ws.Range("E48") "  -4.25%  "
Set-TextValue $ws.Range("D49") "36.70"
Set-TextValue $ws.Range("E49") "  -3.19%  "
Set-TextValue $ws.Range("D50") "0.4193"
Set-TextValue $ws.Range("E50") "  -3.95%  "
Set-TextValue $ws.Range("D51") "0.06036"
Set-TextValue $ws.Range("E51") "  -0.55%  "

Write-Host "Updated cryptos list"
